# Append the 2022-03-02 CRM tank titration result as a new row (61) at the
# bottom of the CRMAccuracyData sheet, and leave the selection where the
# user clicked after entering the data (mirrors the author's "completed
# tank titrations 0302" commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 61

$ws.Cells.Item($row, 1).Value = 20220302
$ws.Cells.Item($row, 2).Value = 2224.779
$ws.Cells.Item($row, 3).Value = 2224.4699999999998
$ws.Cells.Item($row, 4).Formula = "=100*(B$row-C$row)/C$row"
$ws.Cells.Item($row, 5).Value = 180
$ws.Cells.Item($row, 6).Value = "CRM OPENED 20220302"

# Selection moved from I60 to G60 in the saved file.
[void]$ws.Range("G60").Select()
